$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1865
$ws1.Range("F3").Value = 409
$ws1.Range("F5").Value = 868
$ws1.Range("F6").Value = 390
$ws1.Range("F7").Value = 754
$ws1.Range("F8").Value = 13264
$ws1.Range("F9").Value = 13152
$ws1.Range("F13").Value = 553
$ws1.Range("F14").Value = 70
$ws1.Range("F15").Value = 660
$ws1.Range("F16").Value = 2084
$ws1.Range("F17").Value = 63
$ws1.Range("F18").Value = 41
$ws1.Range("F21").Value = 228
$ws1.Range("F22").Value = 280
$ws1.Range("F23").Value = 748
$ws1.Range("F24").Value = 8

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 17
$ws2.Range("F7").Value = 112
$ws2.Range("F8").Value = 11
$ws2.Range("F9").Value = 22

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 27

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1865
$ws4.Range("F4").Value = 409
$ws4.Range("F6").Value = 868
$ws4.Range("F7").Value = 390
$ws4.Range("F9").Value = 754
$ws4.Range("F10").Value = 13264
$ws4.Range("F11").Value = 13152
$ws4.Range("F15").Value = 553
$ws4.Range("F16").Value = 70
$ws4.Range("F17").Value = 660
$ws4.Range("F19").Value = 17
$ws4.Range("F20").Value = 2084
$ws4.Range("F21").Value = 63
$ws4.Range("F22").Value = 41
$ws4.Range("F27").Value = 27
$ws4.Range("F28").Value = 228
$ws4.Range("F29").Value = 280
$ws4.Range("F30").Value = 748
$ws4.Range("F31").Value = 112
$ws4.Range("F32").Value = 11
$ws4.Range("F33").Value = 8
$ws4.Range("F34").Value = 22
